# Fix Training Data Issue (#48)
# The "Date" column (BF) held the wrong date string ("6-11-2012-13") for
# every team row; the NBA stats for this file were actually for 2013-06-11,
# so correct every data row (2-31) in column BF to the right date string.
#
# NumberFormat is forced to Text ("@") before the assignment so that Excel
# stores the value as the literal string "2013-06-11" instead of silently
# re-interpreting it as a date serial number; the format is cleared again
# right after so the cells keep their original (default) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "6-11-2012-13") {
        $cell.NumberFormat = "@"
        $cell.Value = "2013-06-11"
        $cell.ClearFormats()
    }
}
